$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2226943333333333
$ws.Range("N2").Value = 0.668083
$ws.Range("O2").Value = 0.0947101322715019
$ws.Range("P2").Value = 0.09471013227150192
$ws.Range("Q2").Value = 8.227300585635444
$ws.Range("R2").Value = 74.045705270719
$ws.Range("S2").Value = 0.07888373745488642
$ws.Range("T2").Value = 0.07888373745488643

# Row 3
$ws.Range("G3").Value = 36.94436433333333
$ws.Range("H3").Value = 110.833093
$ws.Range("I3").Value = 0.8328964975864823
$ws.Range("J3").Value = 0.8328964975864824
$ws.Range("O3").Value = 0.3538103900551972
$ws.Range("P3").Value = 0.3538103900551972
$ws.Range("Q3").Value = 30.73487872406777
$ws.Range("R3").Value = 276.61390851661
$ws.Range("S3").Value = 0.2946874346866809
$ws.Range("T3").Value = 0.294687434686681

# Row 4
$ws.Range("G4").Value = 36.94436433333333
$ws.Range("H4").Value = 110.833093
$ws.Range("I4").Value = 0.8328964975864823
$ws.Range("J4").Value = 0.8328964975864824
$ws.Range("M4").Value = 1.296707666666667
$ws.Range("N4").Value = 3.890123
$ws.Range("O4").Value = 0.5514794776733007
$ws.Range("P4").Value = 0.5514794776733009
$ws.Range("Q4").Value = 47.90604047115988
$ws.Range("R4").Value = 431.154364240439
$ws.Range("S4").Value = 0.4593253254449149
$ws.Range("T4").Value = 0.459325325444915

# Row 5
$ws.Range("I5").Value = 0.07608399754092349
$ws.Range("J5").Value = 0.07608399754092349
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2226943333333333
$ws.Range("N5").Value = 0.668083
$ws.Range("O5").Value = 0.0947101322715019
$ws.Range("P5").Value = 0.09471013227150192
$ws.Range("Q5").Value = 0.7515530673256666
$ws.Range("R5").Value = 6.763977605930999
$ws.Range("S5").Value = 0.00720592547084549
$ws.Range("T5").Value = 0.00720592547084549

# Row 6
$ws.Range("I6").Value = 0.07608399754092349
$ws.Range("J6").Value = 0.07608399754092349
$ws.Range("O6").Value = 0.3538103900551972
$ws.Range("P6").Value = 0.3538103900551972
$ws.Range("S6").Value = 0.0269193088469128
$ws.Range("T6").Value = 0.0269193088469128

# Row 7
$ws.Range("I7").Value = 0.07608399754092349
$ws.Range("J7").Value = 0.07608399754092349
$ws.Range("M7").Value = 1.296707666666667
$ws.Range("N7").Value = 3.890123
$ws.Range("O7").Value = 0.5514794776733007
$ws.Range("P7").Value = 0.5514794776733009
$ws.Range("Q7").Value = 4.376153670912333
$ws.Range("R7").Value = 39.385383038211
$ws.Range("S7").Value = 0.04195876322316518
$ws.Range("T7").Value = 0.04195876322316519

# Row 8
$ws.Range("G8").Value = 4.037305666666668
$ws.Range("H8").Value = 12.111917
$ws.Range("I8").Value = 0.09101950487259411
$ws.Range("J8").Value = 0.09101950487259411
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2226943333333333
$ws.Range("N8").Value = 0.668083
$ws.Range("O8").Value = 0.0947101322715019
$ws.Range("P8").Value = 0.09471013227150192
$ws.Range("Q8").Value = 0.8990850939012224
$ws.Range("R8").Value = 8.091765845111
$ws.Range("S8").Value = 0.00862046934577
$ws.Range("T8").Value = 0.008620469345770002

# Row 9
$ws.Range("G9").Value = 4.037305666666668
$ws.Range("H9").Value = 12.111917
$ws.Range("I9").Value = 0.09101950487259411
$ws.Range("J9").Value = 0.09101950487259411
$ws.Range("O9").Value = 0.3538103900551972
$ws.Range("P9").Value = 0.3538103900551972
$ws.Range("Q9").Value = 3.358728787898889
$ws.Range("R9").Value = 30.22855909109
$ws.Range("S9").Value = 0.03220364652160344
$ws.Range("T9").Value = 0.03220364652160344

# Row 10
$ws.Range("G10").Value = 4.037305666666668
$ws.Range("H10").Value = 12.111917
$ws.Range("I10").Value = 0.09101950487259411
$ws.Range("J10").Value = 0.09101950487259411
$ws.Range("M10").Value = 1.296707666666667
$ws.Range("N10").Value = 3.890123
$ws.Range("O10").Value = 0.5514794776733007
$ws.Range("P10").Value = 0.5514794776733009
$ws.Range("Q10").Value = 5.235205210643445
$ws.Range("R10").Value = 47.11684689579101
$ws.Range("S10").Value = 0.05019538900522065
$ws.Range("T10").Value = 0.05019538900522066

Write-Output "applied 100 cell updates"